$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

$lastRow = 376
$template = $ws.Range("A376:C376")

function Set-Row {
    param(
        [int]$RowNum,
        [string]$Label,
        [string]$Translation
    )
    $dst = $ws.Range("A" + $RowNum + ":C" + $RowNum)
    [void]$template.Copy($dst)
    $ws.Cells.Item($RowNum, 1).Value = "cs"
    $ws.Cells.Item($RowNum, 2).Value = $Label
    $ws.Cells.Item($RowNum, 3).Value = $Translation
}


Set-Row 377 "lab.setup.menu" "Setupy"
Set-Row 378 "lab.setup.title" "Setupy"
Set-Row 379 "lab.setup.subtitle" "Setup je složení fyzických zařízení použitých pro vapování."
Set-Row 380 "lab.setup.button.create" "Nový setup"
Set-Row 381 "lab.setup.button.list" "Seznam setupů"
Set-Row 382 "lab.setup.create.title" "Nový setup"
Set-Row 383 "lab.setup.create.subtitle" "Setup je poslední součást, která je potřebná pro sledování požitků z vapingu."
Set-Row 384 "lab.setup.create.submit" "Vytvořit setup"
Set-Row 385 "lab.setup.name.label" "Název setupu"
Set-Row 386 "lab.setup.description.label" "Popis"
Set-Row 387 "lab.setup.driptipId.label" "Náústek"
Set-Row 388 "lab.setup.buildId.label" "Build"
Set-Row 389 "lab.setup.modId.label" "Mod"
Set-Row 390 "lab.driptip.tooltip.create" "Vytvořit náústek"
Set-Row 391 "lab.driptip.create.title" "Nový náústek"
Set-Row 392 "lab.driptip.create.subtitle" "Nezdá se to, ale náústky jsou také důležité."
Set-Row 393 "lab.driptip.code.label" "Kód"
Set-Row 394 "lab.driptip.vendorId.label" "Výrobce"
Set-Row 395 "lab.driptip.create.submit" "Vytvořit náústek"
Set-Row 396 "lab.driptip.created.message" "Náústek [{{data.code}}] byl uložen."
Set-Row 397 "lab.build.tooltip.create" "Vytvořit build"
Set-Row 398 "lab.mod.tooltip.create" "Vytvořit mod"
Set-Row 399 "lab.mod.create.title" "Nový mod"
Set-Row 400 "lab.mod.create.subtitle" "Mod obecně zastupuje zařízení, ze kterého lze vapovat."
Set-Row 401 "lab.mod.name.label" "Název modu"
Set-Row 402 "lab.mod.power.label" "Výkon (watty)"
Set-Row 403 "lab.mod.vendorId.label" "Výrobce"
Set-Row 404 "error.Duplicate entry [z_setup_name_unique] of [z_setup]." "Jméno tohoto setupu je již obsazené, použijte prosím jiné."
Set-Row 405 "lab.setup.created.message" "Setup [{{data.name}}] byl uložen."
Set-Row 406 "lab.setup.list.title" "Seznam setupů"
Set-Row 407 "lab.setup.table.name" "Název"
Set-Row 408 "lab.setup.table.driptip" "Náústek"
Set-Row 409 "lab.setup.table.build" "Build"
Set-Row 410 "lab.setup.table.mod" "Mod"
Set-Row 411 "lab.build.inline.atomizer.tooltip" "Atomizér"
Set-Row 412 "lab.build.inline.wraps.tooltip" "Počet otoček na spirálce"


[void]$ws.Range("B401").Select()
Write-Output "Added rows 377-412"
